# Add a new "TeamNum" column (R) to the teams sheet, with each team's
# assigned team-icon/CBS link number, per the commit:
# "Everything now links off the little team icons to CBS"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Cells.Item(1, 18).Value = "TeamNum"

# Row -> TeamNum values (rows 2..13, one per team, matching sheet order)
$teamNums = @{
    2  = 13
    3  = 17
    4  = 10
    5  = 18
    6  = 1
    7  = 14
    8  = 12
    9  = 4
    10 = 5
    11 = 15
    12 = 2
    13 = 16
}

foreach ($row in $teamNums.Keys) {
    $ws.Cells.Item($row, 18).Value = $teamNums[$row]
}

# Reflect the scrolled/selected state from the edit: user scrolled right to
# reveal the new column and left the new cell selected.
$ws.Range("R2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
